# "fix for category for maps" -----------------------------------------
# Rename the sheets, insert two new "map category" rows near the
# existing "农业气候资源" block, and append a brand-new category block
# ("中国行政区划沿革" with its regional subcategories) at the bottom of
# the maps sheet.

$wb = $excel.ActiveWorkbook

$wsMaps = $wb.Worksheets.Item(1)
$wsDocs = $wb.Worksheets.Item(2)
$wsBooks = $wb.Worksheets.Item(3)

# ---- rename sheets ----------------------------------------------------
$wsMaps.Name = "maps"
$wsDocs.Name = "documents"
$wsBooks.Name = "books"

# ---- insert two rows for the new "约翰斯顿皇家现代地图集" / "1933年民国时期地图"
#      sub-category entries right above the existing t37 (农业气候资源) block ----
$wsMaps.Rows("49:50").Insert()

# ---- pre-format all of the new cell ranges (wrapped text, default font) ----
$wsMaps.Range("B49:D50").ClearFormats()
$wsMaps.Range("B49:D50").WrapText = $true

$wsMaps.Range("B64:D64").ClearFormats()
$wsMaps.Range("B64:D64").WrapText = $true

$wsMaps.Range("B65:E70").ClearFormats()
$wsMaps.Range("B65:E70").WrapText = $true

# ---- fill in the new cell values --------------------------------------
$wsMaps.Range("A63").Value = "t24"
$wsMaps.Range("B64").Value = "t2401"

$wsMaps.Range("C49").Value = "gyhsd"
$wsMaps.Range("D49").Value = "约翰斯顿皇家现代地图集"
$wsMaps.Range("B49").Value = "t2205"

$wsMaps.Range("C50").Value = "gmap1933"
$wsMaps.Range("D50").Value = "1933年民国时期地图"
$wsMaps.Range("B50").Value = "t2207"

$wsMaps.Range("C65").Value = "gyg_db"
$wsMaps.Range("D65").Value = "东北地区"
$wsMaps.Range("C66").Value = "gyg_hn"
$wsMaps.Range("D66").Value = "华南地区"
$wsMaps.Range("C67").Value = "gyg_hz"
$wsMaps.Range("D67").Value = "华中地区"
$wsMaps.Range("C68").Value = "gyg_hb"
$wsMaps.Range("D68").Value = "华北地区"
$wsMaps.Range("C69").Value = "gyg_xb"
$wsMaps.Range("D69").Value = "西北地区"
$wsMaps.Range("C70").Value = "gyg_xn"
$wsMaps.Range("D70").Value = "西南地区"

$wsMaps.Range("B65").Value = "t2407"
$wsMaps.Range("B66").Value = "t2402"
$wsMaps.Range("B67").Value = "t2403"
$wsMaps.Range("B68").Value = "t2404"
$wsMaps.Range("B69").Value = "t2405"
$wsMaps.Range("B70").Value = "t2406"

$wsMaps.Range("D64").Value = "华东地区"
$wsMaps.Range("C64").Value = "gyg_hd"

$wsMaps.Range("C63").Value = "china_yg"
$wsMaps.Range("D63").Value = "中国行政区划沿革"

# ---- view state: maps sheet scrolled to G57, books sheet now the
#      selected/active tab (was maps before) ----------------------------
$wsMaps.Activate()
$wsMaps.Range("G57").Select()

$wsBooks.Activate()
$wsBooks.Range("K23").Select()
